# fix(publipostage): Refactor synthetic array /3
#
# Replace the colored-square emoji markers with book emoji, and rename
# the "noir" color label to "bleu":
#   🟥 -> 📕
#   ⬛ -> 📘
#   🟩 -> 📗
#   🟧 -> 📙
#   noir -> bleu

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "🟥"   = "📕"
    "⬛"   = "📘"
    "🟩"   = "📗"
    "🟧"   = "📙"
    "noir" = "bleu"
}

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $used.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($null -ne $val -and $map.ContainsKey($val)) {
            $cell.Value = $map[$val]
        }
    }
}
